$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 22 (R_auto): value/unit were shifted one column too far right
# Before: A22=R_auto, C22=0.8, D22=unitless
# After:  A22=R_auto, B22=0.8, C22=unitless
$ws.Range("B22").Value = 0.8
$ws.Range("C22").Value = "unitless"
$ws.Range("D22").ClearContents()

# --- Remove the blank spacer row 23, shifting rows 24-25 up to 23-24
$ws.Rows("23:23").Delete() | Out-Null

# --- Keep a trailing (now blank) row 25 present, matching the source row
# that used to hold the #Sedimentation/ObservedMAR_oc block before the shift
$ws.Range("A25").Font.Bold = $false
$ws.Rows("25:25").RowHeight = 14.45

# --- Update the saved selection to match the new active cell
$ws.Range("E22").Select() | Out-Null
